# RF011 - Gerenciar Competencias (Portfolio): gender-agreement fixes
# ("o/os/selecionado/cadastrados/excluido" -> "a/as/selecionada/cadastradas/excluida"),
# matching the "1.0 -> 1.1" commit.

$d = $word.ActiveDocument

function Replace-Once($find, $replace) {
    # Search from a fresh Range bound to the whole document each call, so the
    # Nth call naturally lands on the Nth remaining (i.e. first not-yet-fixed)
    # occurrence of $find. Find.Execute collapses $rng to the matched text;
    # assigning .Text (rather than using Find's own Replace argument) swaps
    # in the new text without Word's AutoFormat smart-quote substitution
    # mangling the straight apostrophes used in the source document.
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Replace-Once: text not found: $find"
    }
    $rng.Text = $replace
}

function Replace-Many($find, $replace, $count) {
    for ($i = 0; $i -lt $count; $i++) {
        Replace-Once $find $replace
    }
}

# 1. Basic Flow, step 2
Replace-Once `
    "2. System exibe a listagem dos Competencias (Portfolio) cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' " `
    "2. System exibe a listagem das Competencias (Portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' "

# 2. "System destaca o/a Competencia(s) (Portfolio) selecionado/a na listagem" - appears in AF[1], AF[2], AF[3]
Replace-Many `
    "2. System destaca o Competencias (Portfolio) selecionado na listagem " `
    "2. System destaca a Competencia (Portfolio) selecionada na listagem " `
    3

# 3. AF[1], step 3
Replace-Once `
    "3. Lider de Pessoas clica na opcao 'Editar' para modificar o Competencias (Portfolio) selecionado bs 4" `
    "3. Lider de Pessoas clica na opcao 'Editar' para modificar a Competencia (Portfolio) selecionada bs 4"

# 4. AF[2] heading
Replace-Once `
    "AF[2] – Confirmar Exclusao do Competencias (Portfolio)" `
    "AF[2] – Confirmar Exclusao da Competencia (Portfolio)"

# 5. "Lider de Pessoas clica na opcao 'Excluir' para excluir o/a Competencia(s) (Portfolio) selecionado/a" - AF[2] and AF[3]
Replace-Many `
    "3. Lider de Pessoas clica na opcao 'Excluir' para excluir o Competencias (Portfolio) selecionado " `
    "3. Lider de Pessoas clica na opcao 'Excluir' para excluir a Competencia (Portfolio) selecionada " `
    2

# 6. AF[2], step 6
Replace-Once `
    "6. System exibe a listagem dos Competencias (Portfolio) sem o Competencias (Portfolio) excluido ef[3]" `
    "6. System exibe a listagem das Competencias (Portfolio) sem a Competencia (Portfolio) excluida ef[3]"

# 7. AF[3] heading
Replace-Once `
    "AF[3] – Negar Exclusao do Competencias (Portfolio)" `
    "AF[3] – Negar Exclusao da Competencia (Portfolio)"

# 8. AF[3], step 6
Replace-Once `
    "6. System exibe a listagem dos Competencias (Portfolio) com o Competencias (Portfolio) excluido " `
    "6. System exibe a listagem das Competencias (Portfolio) com a Competencia (Portfolio) excluida "

# 9. AF[4], step 2
Replace-Once `
    "2. System exibe a listagem dos Competencias (Portfolio) cadastrados apenas para visualizacao com a opcao 'Ajuda' " `
    "2. System exibe a listagem das Competencias (Portfolio) cadastradas apenas para visualizacao com a opcao 'Ajuda' "
